# Fix Training Data Issue (#48)
# The "Date" column (BF) held values in the wrong format/day due to how
# NBA stats were originally scraped (e.g. "4-28-2011-12"). Correct this
# to the proper ISO-like date string "2012-04-28" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateColumn = 58   # column BF
$firstRow = 2
$lastRow = 31
$correctDate = "2012-04-28"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateColumn)
    # Assign as literal text (not as an auto-converted date serial number)
    $cell.Value = "'" + $correctDate
    # Restore the default "Normal" style so no extra number-format/style
    # gets attached to the cell (matches original unstyled cells).
    $cell.Style = "Normal"
}
